$d = $word.ActiveDocument

function Set-ParaText($paraIndex, $newText) {
    $p = $d.Paragraphs.Item($paraIndex)
    $rng = $p.Range
    $target = $d.Range($rng.Start, $rng.End - 1)
    $target.Text = $newText
}

function Replace-InPara($paraIndex, $oldText, $newText) {
    $p = $d.Paragraphs.Item($paraIndex)
    $rng = $p.Range
    $null = $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
}

# ============================================================
# Work from the BOTTOM of the document upward so that paragraph
# indices for not-yet-processed (earlier) paragraphs stay valid.
# Paragraph numbers below refer to the ORIGINAL document.
# ============================================================

# --- Paragraphs 16 (empty) and 17 (last) are unchanged. ---

# --- Paragraph 15: "As I kept on exploring ... this course." ---
# Rewritten & extended; also a brand-new paragraph is added right after it.
Set-ParaText 15 "As I kept on exploring the different ideas and algorithms used in the above mentioned modules, the exposure increased my inclination towards learning more about these analysis techniques and therefore, a formalised classroom training process would be the ideal way gain further knowledge of such an insightful subject. Growing with this point of view has finally brought me to this juncture where I have all the required pre-requisites, motivation and goals to take this up as a career and finally contribute to some……. [Current college name] is reputed in imparting unparalleled knowledge about data science and hence will be the right choice possible to fulfil my academic aspirations."

# Make "Current college name" bold within paragraph 15.
$p15 = $d.Paragraphs.Item(15)
$boldRng = $p15.Range.Duplicate()
$found = $boldRng.Find.Execute("Current college name")
if ($found) {
    $boldRng.Bold = 1
}

# Insert the new paragraph "What I have to offer to the society." after paragraph 15.
$p15 = $d.Paragraphs.Item(15)
$p15.Range.InsertParagraphAfter()
Set-ParaText 16 "What I have to offer to the society."

# --- Paragraph 14: "This is how I developed ... Pandas." ---
Replace-InPara 14 "I moved towards modules for large data handling" "I moved towards python modules for large data handling"

# --- Paragraph 13 (empty) gets removed entirely. ---
$d.Paragraphs.Item(13).Range.Delete()

# --- Paragraph 12 ("Add a tinge of imformalness.") list item removed entirely. ---
$d.Paragraphs.Item(12).Range.Delete()

# --- Paragraph 11 ("Root cause analysis.") list item removed entirely. ---
$d.Paragraphs.Item(11).Range.Delete()

# --- Paragraph 10: "Working in the financial domain ... eliminated." ---
Set-ParaText 10 "Working in the financial domain as my project in TCS I have handled large amount of critical data. And at the same time being in a support project, I noticed that there is a lot of manual analysis that was both time consuming and redundant. There were a lot of unknowns and hit and trials involved. And all this was because we did not have a root cause analysis model in place. This encouraged me to explore the different ways in which we can have a better idea about the problem before we can jump to the solution."

# --- Paragraph 9: "TCS gave me the opportunity ... critical files." gets merged
#     into paragraph 8, after paragraph 8's own text is updated. ---

# First remove the (now relocated) _GoBack bookmark that currently sits inside
# paragraph 9, so it does not end up in the wrong place.
$d.Bookmarks.ShowHidden = $true
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

Set-ParaText 9 " I joined TCS as an associate trainee where I made use of these learnings in real world problem solving. I started with automating simple day to day tasks like checking the quality of files delivered to us based on company standards, which critical files are still left to be delivered and creating a daily log file containing all the manual changes that was made to certain critical files. "

# Merge paragraph 9 into paragraph 8 (delete the paragraph mark between them).
$p8 = $d.Paragraphs.Item(8)
$markRng = $d.Range($p8.Range.End - 1, $p8.Range.End)
$markRng.Delete()

# --- Paragraph 8: update its own tail text (now paragraph 8, still). ---
Replace-InPara 8 "I had a lot of interest in these courses and completed these with 5-star ratings." "After completing these courses with 5-star ratings"

# --- Paragraph 7: "Programming languages ... semester." -> "... year." ---
Replace-InPara 7 " of the semester.", " of the year."

# --- Paragraph 6 (empty) stays; a brand-new paragraph is inserted right after it. ---
$p6 = $d.Paragraphs.Item(6)
$p6.Range.InsertParagraphAfter()
Set-ParaText 7 "College opened my ways to learn different technologies. "

# --- Paragraph 5: "College opened my ways to learn different technologies. My
#     affinity towards technology…." is removed entirely (its content has been
#     redistributed: the new paragraph above, and a tail appended to paragraph 4). ---
$d.Paragraphs.Item(5).Range.Delete()

# --- Paragraph 4: append new tail text. ---
Replace-InPara 4 "not to miss out on those opportunities. ", "not to miss out on those opportunities. And my affinity towards technology have always helped me"

# --- Paragraphs 2 and 3 ("Write analysis as a whole" / "….think of a reason …") unchanged. ---

# --- Paragraph 1: "Computers has been very intriguing ..." completely rewritten
#     and split into two new paragraphs. ---
Set-ParaText 1 "The subject of computers have always intrigued me. So, whether it be learning GW Basic or LOGO in elementary school or learning algorithms and data structures during my under-graduation, I have always learnt them with utmost enthusiasm. Given a computational problem and exploring different ways and means to solve them has always interested me more than anything else."

$p1 = $d.Paragraphs.Item(1)
$p1.Range.InsertParagraphAfter()
Set-ParaText 2 "It has always occurred to me that, in any kind of problem, the ‘HOW’ part in a problem is always important than the final solution itself. Proper insights about the problem like, identifying cause and patterns, possible solution alternatives and their pros and cons, will always help in providing a robust solution. This is the reason, given any problem, I have always spent more time analysing and planning rather developing and implementation itself. And this is the same strategy that I have always had in mind since I wrote my first computer program. "

# Re-create the _GoBack bookmark roughly where the last edit happened (inside the
# newly written first paragraph), matching the relocation implied by the diff.
$p1 = $d.Paragraphs.Item(1)
$bmRng = $p1.Range.Duplicate()
$null = $bmRng.Find.Execute("always intrigued me")
if ($bmRng.Find.Found) {
    $d.Bookmarks.Add("_GoBack", $bmRng)
}

Write-Host "Done. Final paragraph count:" $d.Paragraphs.Count
